$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows right before the old "GRAND TOTAL" row (old row 20),
# which become rows 20-24. Excel carries the surrounding banded
# fill/border formatting into the freshly inserted rows automatically,
# and everything below (old rows 20-39) shifts down to rows 25-44.
$ws.Rows("20:24").Insert()

# Populate the new "flow" rows with the specification text added in this
# commit (exchange add/remove/dropdown/reset/placeholder flow notes).
$ws.Range("C20").Value = "Add exchange button"
$ws.Range("C21").Value = "Remove exchange2 row"
$ws.Range("C22").Value = "Add dropdown and populate with added exchanges"
$ws.Range("C23").Value = "(reset button)"
$ws.Range("C24").Value = "(placeholder text)"

# Match the author's final selection/viewport noted in the workbook.
$ws.Range("D31").Select()
